$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 4.805799999999998
$ws.Range("A8").Value = -22.42270000000002
$ws.Range("A10").Value = -21.75659999999998
$ws.Range("A12").Value = -21.56539999999999
$ws.Range("B15").Value = 4.752499999999995
$ws.Range("A18").Value = -21.78590000000002
$ws.Range("B18").Value = 5.976299999999997
$ws.Range("B20").Value = 9.372599999999991
$ws.Range("B29").Value = 4.850300000000004
$ws.Range("B30").Value = 5.665599999999999
$ws.Range("B31").Value = 5.376100000000001
$ws.Range("A37").Value = -19.5289
$ws.Range("B40").Value = 9.312299999999995
$ws.Range("B50").Value = 5.1114
$ws.Range("A55").Value = -21.9229
$ws.Range("A68").Value = -21.485
$ws.Range("B68").Value = 4.5547
$ws.Range("B76").Value = 6.033599999999996
$ws.Range("A77").Value = -20.43079999999999
$ws.Range("A78").Value = -19.88019999999999
$ws.Range("A81").Value = -21.8081
$ws.Range("A82").Value = -21.86259999999999
$ws.Range("B87").Value = 4.795499999999995
$ws.Range("B88").Value = 4.414399999999998
$ws.Range("B96").Value = 5.313800000000006
$ws.Range("B98").Value = 6.513499999999999
$ws.Range("B101").Value = 9.639199999999995
$ws.Range("B102").Value = 8.226500000000007
